$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.937.39"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "3.505.74"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("D7").Value = "3.503.07"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.23%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "4.095.90"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "67.928.85"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "3.496.80"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").Value = "3.643.89"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -3.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "3.499.93"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "176.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0900"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "30.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.83%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.897"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "
